# Implement followed and ctx eligible
# Fills in row 22 ("CTX_ELIGIBLE" breakdown header) on the "File active" sheet,
# mirroring the existing row 21 ("ACTIVE_LIST") layout/format, bumps the
# height of rows 22 and 26, and moves the sheet's selection down to A26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 22: copy formatting from row 21 (B:AB) --------------------------
$ws.Range("B21:AB21").Copy()
$ws.Range("B22:AB22").PasteSpecial(-4122)

# --- Row 22: fill in the CTX_ELIGIBLE breakdown values --------------------
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
foreach ($col in $cols) {
    $srcAddr = $col + "21"
    $dstAddr = $col + "22"
    $srcText = $ws.Range($srcAddr).Value2
    $dstText = $srcText -replace "ACTIVE_LIST", "CTX_ELIGIBLE"
    $ws.Range($dstAddr).Value = $dstText
}

# --- Row heights -----------------------------------------------------------
$ws.Rows.Item(22).RowHeight = 34.55
$ws.Rows.Item(26).RowHeight = 28.25

# --- Selection / view: move to A26 -----------------------------------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $True
$ws.Activate() | Out-Null
$ws.Range("A26").Select() | Out-Null
$win.ScrollRow = 21
